# Add a new row to the "Completed" reading list for "A Marvelous Life"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$row = 17

$ws.Cells.Item($row, 1).Value = "A Marvelous Life"
$ws.Cells.Item($row, 2).Value = "Danny Fingeroth"

$ws.Cells.Item($row, 3).Value = 43853
$ws.Cells.Item($row, 4).Value = 43859

# Reuse the exact date formatting (style) already used by the column above
$ws.Range("C16:D16").Copy()
$ws.Range("C17:D17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item($row, 5).Value = "biography;marvel;stan lee;comics"
$ws.Cells.Item($row, 6).Value = "Audio"
$ws.Cells.Item($row, 7).Value = "14 Hrs 47 Mins"
